# Weekly update: insert a new record at the top of the data table (row 59),
# pushing the existing rows 59-69 down to 60-70, and populate the new row
# with this week's price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(59).Insert()

$ws.Cells.Item(59, 1).Value  = 11
$ws.Cells.Item(59, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(59, 3).Value  = "Bíobío"
$ws.Cells.Item(59, 4).Value  = 44951
$ws.Cells.Item(59, 5).Value  = 8
$ws.Cells.Item(59, 6).Value  = 100112031
$ws.Cells.Item(59, 7).Value  = "Poroto verde"
$ws.Cells.Item(59, 8).Value  = "Magnum"
$ws.Cells.Item(59, 9).Value  = "Primera"
$ws.Cells.Item(59, 10).Value = 100
$ws.Cells.Item(59, 11).Value = 22000
$ws.Cells.Item(59, 12).Value = 24000
$ws.Cells.Item(59, 13).Value = 23000
$ws.Cells.Item(59, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(59, 15).Value = "Región Metropolitana"
$ws.Cells.Item(59, 16).Value = 920
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"
